$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 83 ("「宇宙の奇石」..." post) entirely; all subsequent rows shift up by one.
$ws.Rows.Item(83).Delete()
